# Update "想去人数" (interested-count) values in column F for a handful of
# events that appear on both the "展览" sheet and the aggregated "全部类型"
# sheet. The row numbers differ between the two sheets because "全部类型"
# interleaves events from the other category sheets, so each sheet is
# updated using its own row indices.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (worksheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 181
$ws1.Range("F12").Value = 5246
$ws1.Range("F13").Value = 67
$ws1.Range("F14").Value = 865
$ws1.Range("F15").Value = 125
$ws1.Range("F16").Value = 2325
$ws1.Range("F19").Value = 2170

# --- Sheet "全部类型" (worksheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 181
$ws4.Range("F12").Value = 5246
$ws4.Range("F14").Value = 67
$ws4.Range("F16").Value = 865
$ws4.Range("F17").Value = 125
$ws4.Range("F18").Value = 2325
$ws4.Range("F22").Value = 2170
